$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.871.38'
$ws.Range("E2").Value = '  +0.63%  '

$ws.Range("D3").Value = '1.642.39'
$ws.Range("E3").Value = '  +0.56%  '

$ws.Range("E4").Value = '  -0.78%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.76%  '

$ws.Range("E6").Value = '  +2.16%  '

$ws.Range("E7").Value = '  -0.74%  '

$ws.Range("E8").Value = '  +1.96%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0622'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.34%  '

$ws.Range("D12").Value = '1.871.68'
$ws.Range("E12").Value = '  +0.50%  '

$ws.Range("D13").Value = '1.639.29'
$ws.Range("E13").Value = '  +0.51%  '

$ws.Range("E14").Value = '  +0.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.84%  '

$ws.Range("D17").Value = '26.876.32'
$ws.Range("E17").Value = '  +0.70%  '

$ws.Range("E18").Value = '  +1.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '219.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.88%  '

$ws.Range("E20").Value = '  -0.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.66'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.08%  '

$ws.Range("E23").Value = '  +3.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.32%  '

$ws.Range("E26").Value = '  -0.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.52%  '

$ws.Range("E28").Value = '  +1.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.82'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0511'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.11%  '

$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("E32").Value = '  +0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.99'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.13%  '

$ws.Range("E34").Value = '  +3.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").Value = '1.240.57'
$ws.Range("E36").Value = '  -1.38%  '

$ws.Range("E37").Value = '  +1.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.540'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.836'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.80%  '

$ws.Range("E40").Value = '  -0.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.807'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.00%  '

$ws.Range("D43").Value = '1.782.69'
$ws.Range("E43").Value = '  +0.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '60.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.66%  '

$ws.Range("E47").Value = '  +0.82%  '

$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  +13.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0974'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.01%  '
